$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 19: H19,J19,L19,N19
$ws.Range("H19").Value = 908.7895
$ws.Range("J19").Value = 983.53845
$ws.Range("L19").Value = 983.53845
$ws.Range("N19").Value = -1333.53845
# Row 39: H39,I39,J39,K39,L39,M39,N39
$ws.Range("H39").Value = 883.86664
$ws.Range("I39").Value = 1330.8889
$ws.Range("J39").Value = 213.33333
$ws.Range("K39").Value = 3992.6667
$ws.Range("L39").Value = 639.99999
$ws.Range("M39").Value = -3696.6667
$ws.Range("N39").Value = -1231.99999
# Row 42: H42,I42,J42,K42,L42,M42,N42
$ws.Range("H42").Value = 3750.3333
$ws.Range("I42").Value = 5575.5
$ws.Range("J42").Value = 100
$ws.Range("K42").Value = 16726.5
$ws.Range("L42").Value = 300
$ws.Range("M42").Value = -16496.5
$ws.Range("N42").Value = -760
# Row 106: H106,I106,K106,M106
$ws.Range("H106").Value = 37039536
$ws.Range("I106").Value = 37039536
$ws.Range("K106").Value = 37039536
$ws.Range("M106").Value = -37038905
# Row 129: H129,I129,J129,K129,L129,M129,N129
$ws.Range("H129").Value = 1023.35
$ws.Range("I129").Value = 464.25
$ws.Range("J129").Value = 1085.4722
$ws.Range("K129").Value = 1392.75
$ws.Range("L129").Value = 3256.4166
$ws.Range("M129").Value = 3607.25
$ws.Range("N129").Value = -13256.4166
# Row 137: H137,I137,J137,K137,L137,M137,N137
$ws.Range("H137").Value = 1637.6897
$ws.Range("I137").Value = 1680.381
$ws.Range("J137").Value = 1525.625
$ws.Range("K137").Value = 5041.143
$ws.Range("L137").Value = 4576.875
$ws.Range("M137").Value = -2491.143
$ws.Range("N137").Value = -9676.875
# Row 140: H140,J140,L140,N140
$ws.Range("H140").Value = 45333.332
$ws.Range("J140").Value = 45333.332
$ws.Range("L140").Value = 45333.332
$ws.Range("N140").Value = -55693.332

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 74: H74,I74,J74,K74,L74,M74,N74
$ws.Range("H74").Value = 6165.2173
$ws.Range("I74").Value = 993.6
$ws.Range("J74").Value = 40642.668
$ws.Range("K74").Value = 993.6
$ws.Range("L74").Value = 40642.668
$ws.Range("M74").Value = -119.6
$ws.Range("N74").Value = -42390.668
# Row 77: H77,I77,J77,K77,L77,M77,N77
$ws.Range("H77").Value = 6165.2173
$ws.Range("I77").Value = 993.6
$ws.Range("J77").Value = 40642.668
$ws.Range("K77").Value = 4968
$ws.Range("L77").Value = 203213.34
$ws.Range("M77").Value = -600
$ws.Range("N77").Value = -211949.34
# Row 110: H110,I110,K110,M110
$ws.Range("H110").Value = 1160.45
$ws.Range("I110").Value = 1075.5625
$ws.Range("K110").Value = 1075.5625
$ws.Range("M110").Value = 969.4375
# Row 139: H139,J139,L139,N139
$ws.Range("H139").Value = 64238.332
$ws.Range("J139").Value = 64238.332
$ws.Range("L139").Value = 64238.332
$ws.Range("N139").Value = -74518.33199999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 59: H59,J59,L59,N59
$ws.Range("H59").Value = 46775
$ws.Range("J59").Value = 46775
$ws.Range("L59").Value = 46775
$ws.Range("N59").Value = -48469
# Row 133: H133,J133,L133,N133
$ws.Range("H133").Value = 52000
$ws.Range("J133").Value = 52000
$ws.Range("L133").Value = 52000
$ws.Range("N133").Value = -62120

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31,I31,J31,K31,L31,M31,N31
$ws.Range("H31").Value = 1623.9375
$ws.Range("I31").Value = 1067.9333
$ws.Range("J31").Value = 9964
$ws.Range("K31").Value = 1067.9333
$ws.Range("L31").Value = 9964
$ws.Range("M31").Value = -772.9332999999999
$ws.Range("N31").Value = -10554
# Row 34: H34,I34,J34,K34,L34,M34,N34
$ws.Range("H34").Value = 1623.9375
$ws.Range("I34").Value = 1067.9333
$ws.Range("J34").Value = 9964
$ws.Range("K34").Value = 1067.9333
$ws.Range("L34").Value = 9964
$ws.Range("M34").Value = -865.9332999999999
$ws.Range("N34").Value = -10368

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5: H5,I5,J5,K5,L5,M5,N5
$ws.Range("H5").Value = 2525.3845
$ws.Range("I5").Value = 2207.7144
$ws.Range("J5").Value = 2896
$ws.Range("K5").Value = 6623.1432
$ws.Range("L5").Value = 8688
$ws.Range("M5").Value = -6511.1432
$ws.Range("N5").Value = -8912
# Row 70: H70,I70,J70,K70,L70,M70,N70
$ws.Range("H70").Value = 3199.8572
$ws.Range("I70").Value = 699.5
$ws.Range("J70").Value = 4200
$ws.Range("K70").Value = 2098.5
$ws.Range("L70").Value = 12600
$ws.Range("M70").Value = -1783.5
$ws.Range("N70").Value = -13230
# Row 73: H73,I73,J73,K73,L73,M73,N73
$ws.Range("H73").Value = 3199.8572
$ws.Range("I73").Value = 699.5
$ws.Range("J73").Value = 4200
$ws.Range("K73").Value = 2098.5
$ws.Range("L73").Value = 12600
$ws.Range("M73").Value = -1006.5
$ws.Range("N73").Value = -14784
# Row 75: H75,I75,J75,K75,L75,M75,N75
$ws.Range("H75").Value = 2889
$ws.Range("I75").Value = 2563
$ws.Range("J75").Value = 2982.1428
$ws.Range("K75").Value = 7689
$ws.Range("L75").Value = 8946.428400000001
$ws.Range("M75").Value = -6691
$ws.Range("N75").Value = -10942.4284
# Row 78: H78,I78,J78,K78,L78,M78,N78
$ws.Range("H78").Value = 2889
$ws.Range("I78").Value = 2563
$ws.Range("J78").Value = 2982.1428
$ws.Range("K78").Value = 23067
$ws.Range("L78").Value = 26839.2852
$ws.Range("M78").Value = -18075
$ws.Range("N78").Value = -36823.2852
# Row 121: H121,I121,K121,M121
$ws.Range("H121").Value = 195
$ws.Range("I121").Value = 195
$ws.Range("K121").Value = 585
$ws.Range("M121").Value = 725
# Row 122: H122,J122,L122,N122
$ws.Range("H122").Value = 867.3333
$ws.Range("J122").Value = 1532.6666
$ws.Range("L122").Value = 13793.9994
$ws.Range("N122").Value = -18693.9994
# Row 131: H131,I131,J131,K131,L131,M131,N131
$ws.Range("H131").Value = 2439.9102
$ws.Range("I131").Value = 306
$ws.Range("J131").Value = 2566.9285
$ws.Range("K131").Value = 918
$ws.Range("L131").Value = 7700.7855
$ws.Range("M131").Value = 4122
$ws.Range("N131").Value = -17780.7855
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 1096.1364
$ws.Range("I132").Value = 1023.8889
$ws.Range("J132").Value = 1146.1538
$ws.Range("K132").Value = 9215.000100000001
$ws.Range("L132").Value = 10315.3842
$ws.Range("M132").Value = -6685.000100000001
$ws.Range("N132").Value = -15375.3842
# Row 135: H135,I135,J135,K135,L135,M135,N135
$ws.Range("H135").Value = 2525.3845
$ws.Range("I135").Value = 2207.7144
$ws.Range("J135").Value = 2896
$ws.Range("K135").Value = 19869.4296
$ws.Range("L135").Value = 26064
$ws.Range("M135").Value = -17334.4296
$ws.Range("N135").Value = -31134
# Row 137: H137,I137,J137,K137,L137,M137,N137
$ws.Range("H137").Value = 5319184.5
$ws.Range("I137").Value = 8334789.5
$ws.Range("J137").Value = 149576.14
$ws.Range("K137").Value = 25004368.5
$ws.Range("L137").Value = 448728.42
$ws.Range("M137").Value = -24999268.5
$ws.Range("N137").Value = -458928.42

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 43: H43,J43,L43,N43
$ws.Range("H43").Value = 12000
$ws.Range("J43").Value = 12000
$ws.Range("L43").Value = 12000
$ws.Range("N43").Value = -12302
# Row 107: H107,J107,L107,N107
$ws.Range("H107").Value = 1000.3077
$ws.Range("J107").Value = 478.66666
$ws.Range("L107").Value = 478.66666
$ws.Range("N107").Value = -4318.66666
# Row 138: H138,J138,L138,N138
$ws.Range("H138").Value = 62500
$ws.Range("J138").Value = 62500
$ws.Range("L138").Value = 62500
$ws.Range("N138").Value = -72780

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22: H22,I22,J22,K22,L22,M22,N22
$ws.Range("H22").Value = 6593.4116
$ws.Range("I22").Value = 97.75
$ws.Range("J22").Value = 8592.076999999999
$ws.Range("K22").Value = 97.75
$ws.Range("L22").Value = 8592.076999999999
$ws.Range("M22").Value = 197.25
$ws.Range("N22").Value = -9182.076999999999
# Row 27: H27,I27,J27,K27,L27,M27,N27
$ws.Range("H27").Value = 6593.4116
$ws.Range("I27").Value = 97.75
$ws.Range("J27").Value = 8592.076999999999
$ws.Range("K27").Value = 97.75
$ws.Range("L27").Value = 8592.076999999999
$ws.Range("M27").Value = 9.25
$ws.Range("N27").Value = -8806.076999999999
# Row 46: H46,J46,L46,N46
$ws.Range("H46").Value = 2626.6667
$ws.Range("J46").Value = 2626.6667
$ws.Range("L46").Value = 2626.6667
$ws.Range("N46").Value = -3002.6667
# Row 55: H55,I55,J55,K55,L55,M55,N55
$ws.Range("H55").Value = 209.47826
$ws.Range("I55").Value = 142.29411
$ws.Range("J55").Value = 399.83334
$ws.Range("K55").Value = 142.29411
$ws.Range("L55").Value = 399.83334
$ws.Range("M55").Value = 30.70589000000001
$ws.Range("N55").Value = -745.83334

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 36156
$ws.Range("I122").Value = 59913.47
$ws.Range("J122").Value = 2499.5833
$ws.Range("K122").Value = 179740.41
$ws.Range("L122").Value = 7498.749899999999
$ws.Range("M122").Value = -177290.41
$ws.Range("N122").Value = -12398.7499
